$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows 2-17 to reflect caminhão3 (truck 3) new delivery roster ---
$ws.Range("A2").Value = "'11117"
$ws.Range("B2").Value = "21/11/2024"
$ws.Range("C2").Value = "manhã"
$ws.Range("D2").Value = "Desconhecido"
$ws.Range("E2").Value = "Campinas"

$ws.Range("A3").Value = "'11117"
$ws.Range("B3").Value = "21/11/2024"
$ws.Range("C3").Value = "manhã"
$ws.Range("D3").Value = "Desconhecido"
$ws.Range("E3").Value = "Campinas"

$ws.Range("A4").Value = "'11117"
$ws.Range("B4").Value = "21/11/2024"
$ws.Range("C4").Value = "manhã"
$ws.Range("D4").Value = "Desconhecido"
$ws.Range("E4").Value = "Campinas"

$ws.Range("A5").Value = "'11117"
$ws.Range("B5").Value = "21/11/2024"
$ws.Range("C5").Value = "manhã"
$ws.Range("D5").Value = "Desconhecido"
$ws.Range("E5").Value = "Campinas"

$ws.Range("A6").Value = "'11136"
$ws.Range("B6").Value = "21/11/2024"
$ws.Range("C6").Value = "manhã"
$ws.Range("D6").Value = "Desconhecido"
$ws.Range("E6").Value = "Centro"

$ws.Range("A7").Value = "'11179"
$ws.Range("B7").Value = "21/11/2024"
$ws.Range("C7").Value = "tarde"
$ws.Range("D7").Value = "Desconhecido"
$ws.Range("E7").Value = "Nova Palhoça"

$ws.Range("A8").Value = "'11179"
$ws.Range("B8").Value = "21/11/2024"
$ws.Range("C8").Value = "tarde"
$ws.Range("D8").Value = "Desconhecido"
$ws.Range("E8").Value = "Nova Palhoça"

$ws.Range("A9").Value = "'11179"
$ws.Range("B9").Value = "21/11/2024"
$ws.Range("C9").Value = "tarde"
$ws.Range("D9").Value = "Desconhecido"
$ws.Range("E9").Value = "Nova Palhoça"

$ws.Range("A10").Value = "'11179"
$ws.Range("B10").Value = "21/11/2024"
$ws.Range("C10").Value = "tarde"
$ws.Range("D10").Value = "Desconhecido"
$ws.Range("E10").Value = "Nova Palhoça"

$ws.Range("A11").Value = "'11180"
$ws.Range("B11").Value = "21/11/2024"
$ws.Range("C11").Value = "tarde"
$ws.Range("D11").Value = "Desconhecido"
$ws.Range("E11").Value = "Nova Palhoça"

$ws.Range("A12").Value = "'11180"
$ws.Range("B12").Value = "21/11/2024"
$ws.Range("C12").Value = "tarde"
$ws.Range("D12").Value = "Desconhecido"
$ws.Range("E12").Value = "Nova Palhoça"

$ws.Range("A13").Value = "'11183"
$ws.Range("B13").Value = "22/11/2024"
$ws.Range("C13").Value = "manhã"
$ws.Range("D13").Value = "Desconhecido"
$ws.Range("E13").Value = "Pagani"

$ws.Range("A14").Value = "'11183"
$ws.Range("B14").Value = "22/11/2024"
$ws.Range("C14").Value = "manhã"
$ws.Range("D14").Value = "Desconhecido"
$ws.Range("E14").Value = "Pagani"

$ws.Range("A15").Value = "'11183"
$ws.Range("B15").Value = "22/11/2024"
$ws.Range("C15").Value = "manhã"
$ws.Range("D15").Value = "Desconhecido"
$ws.Range("E15").Value = "Pagani"

$ws.Range("A16").Value = "'11183"
$ws.Range("B16").Value = "22/11/2024"
$ws.Range("C16").Value = "manhã"
$ws.Range("D16").Value = "Desconhecido"
$ws.Range("E16").Value = "Pagani"

$ws.Range("A17").Value = "'11185"
$ws.Range("B17").Value = "22/11/2024"
$ws.Range("C17").Value = "manhã"
$ws.Range("D17").Value = "Desconhecido"
$ws.Range("E17").Value = "Passa Vinte"

# --- Append new rows 18-22; match the center/center alignment style used by existing data rows ---
$newRange = $ws.Range("A18:E22")
$newRange.HorizontalAlignment = -4108  # xlCenter
$newRange.VerticalAlignment = -4108    # xlCenter

$ws.Range("A18").Value = "'11158"
$ws.Range("B18").Value = "22/11/2024"
$ws.Range("C18").Value = "tarde"
$ws.Range("D18").Value = "Desconhecido"
$ws.Range("E18").Value = "Forquilhinha"

$ws.Range("A19").Value = "'11158"
$ws.Range("B19").Value = "22/11/2024"
$ws.Range("C19").Value = "tarde"
$ws.Range("D19").Value = "Desconhecido"
$ws.Range("E19").Value = "Forquilhinha"

$ws.Range("A20").Value = "'11158"
$ws.Range("B20").Value = "22/11/2024"
$ws.Range("C20").Value = "tarde"
$ws.Range("D20").Value = "Desconhecido"
$ws.Range("E20").Value = "Forquilhinha"

$ws.Range("A21").Value = "'11158"
$ws.Range("B21").Value = "22/11/2024"
$ws.Range("C21").Value = "tarde"
$ws.Range("D21").Value = "Desconhecido"
$ws.Range("E21").Value = "Forquilhinha"

$ws.Range("A22").Value = "'11154"
$ws.Range("B22").Value = "22/11/2024"
$ws.Range("C22").Value = "tarde"
$ws.Range("D22").Value = "Desconhecido"
$ws.Range("E22").Value = "Forquilhas"

